$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$newRow = 65
$prevRow = 64

$ws.Range("A" + $prevRow + ":E" + $prevRow).Copy()
$ws.Range("A" + $newRow + ":E" + $newRow).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = 43966
$ws.Cells.Item($newRow, 2).Value = 38565
$ws.Cells.Item($newRow, 3).Value = 1708
$ws.Cells.Item($newRow, 4).Value = 48
$ws.Cells.Item($newRow, 5).Value = 2103

$tbl = $ws.ListObjects.Item("Table3")
$tbl.Resize($ws.Range("A1:E" + $newRow))

$ws.Range("E65").Select()
